$wb = $excel.ActiveWorkbook

# --- Fix "Chasity" -> "Chastity" typo on the Characters sheet ---
$chars = $wb.Worksheets.Item("Characters")
$chars.Range("C2").Value = "Chastity"
$chars.Range("F10").Value = "Chastity"

# --- Update view/selection state ---
# Enemies sheet: selection moves from A9 to C4:D5 (active cell D4), no longer the
# selected (visible) tab.
$enemies = $wb.Worksheets.Item("Enemies")
$enemies.Activate()
$enemies.Range("C4:D5").Select()

# Characters sheet becomes the selected (visible) tab, with selection D11.
$chars.Activate()
$chars.Range("D11").Select()

Write-Host "done"
